# Revert capacity chart to show kilowatts on the y-axis.
#
# - Solar capacity figures (column E / rows 21-26) were stored in watts;
#   convert them back to kilowatts (divide by 1000).
# - The shared number format used by the data table (columns B:G) goes from
#   a plain integer format to one that keeps a single decimal place, so the
#   new kW figures with fractional parts still display correctly.
# - The value (y) axis of the embedded chart is relabeled "Kilowatts (kW)"
#   (was "Watts") and its tick-label number format reverts from the
#   "show thousands as K" custom format to a plain "#,##0".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Worksheet data: convert Solar (column E) watts -> kilowatts ---------
$ws.Range("E21").Value = 13.8
$ws.Range("E22").Value = 11.2
$ws.Range("E23").Value = 25.1
$ws.Range("E24").Value = 50.81
$ws.Range("E25").Value = 79.38
$ws.Range("E26").Value = 45.93

# --- Number format: one decimal place so fractional kW values show up ----
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- Chart: relabel the value axis and fix its tick number format --------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
